$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet tab / workbook.xml <sheet name="..."/>
$ws.Name = "Trail Balance"

# Remove the old journal-report rows, leaving gaps (do not shift remaining rows)
$ws.Rows.Item(6).ClearContents()
$ws.Rows.Item(8).ClearContents()
$ws.Rows.Item(9).ClearContents()
$ws.Rows.Item(12).ClearContents()
$ws.Rows.Item(14).ClearContents()
$ws.Rows.Item(15).ClearContents()

# New report title in A1
$ws.Range("A1").Value = "ElHadar-PLC  Statement of profit or Loss and other comprehensive income  From Fri Jan 01 1999 To Sun Jan 01 2023"

# New header row for the trial-balance-style table
$ws.Range("A4").Value = "Account"
$ws.Range("C4").Value = "Total"

# B4 stays blank (empty string cell, like the source template); write it via a
# quote-prefixed empty value and then drop back to the Normal style so no
# stray quote-prefix formatting is left behind
$ws.Range("B4").Value = "'"
$ws.Range("B4").Style = "Normal"

# Add a 4th column matching the existing column widths
$ws.Columns.Item(4).ColumnWidth = 50

# Grow the title merge block to include the new column
$ws.Range("A1:C3").UnMerge()
$ws.Range("A1:D3").Merge()
